# Replace the static "www.drpaulduenas.com" text in the footer with a
# configurable MERGEFIELD ("=website") so the website can come from the
# mail-merge / template data source instead of being hard-coded.
#
# Target OOXML (inside the existing <w:p w14:paraId="49FEE76D" .../> footer
# paragraph) goes from a single run:
#
#   <w:r><w:rPr>...</w:rPr><w:t>www.drpaulduenas.com</w:t></w:r>
#
# to the classic Word field-code run sequence:
#
#   <w:r><w:rPr>...</w:rPr><w:fldChar w:fldCharType="begin"/></w:r>
#   <w:r><w:rPr>...</w:rPr><w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r>
#   <w:r><w:rPr>...</w:rPr><w:fldChar w:fldCharType="separate"/></w:r>
#   <w:r><w:rPr>...</w:rPr><w:t>«=website»</w:t></w:r>
#   <w:r><w:rPr>...</w:rPr><w:fldChar w:fldCharType="end"/></w:r>
#
# with the run properties (rFonts Avenir Book, bold, sz/szCs 20) preserved
# unchanged on every run.

$d = $word.ActiveDocument

$targetText = "www.drpaulduenas.com"
$found = $false
$hitRange = $null

# The text lives in a footer, so walk every section's footers (primary,
# first-page, even-page) looking for it instead of assuming the body.
foreach ($sec in $d.Sections) {
    foreach ($f in $sec.Footers) {
        if (-not $f.Exists) { continue }

        $rng = $f.Range.Duplicate
        if ($rng.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
            $found = $true
            $hitRange = $rng
            break
        }
    }
    if ($found) { break }
}

if (-not $found) {
    throw "Could not find '$targetText' in any footer"
}

# Remove the old run's text; this collapses $hitRange to an insertion point
# right where the paragraph's lone run used to be, while leaving the
# paragraph mark (and its pPr / paraId / rsid attributes) untouched.
$hitRange.Delete()

# Rebuild the paragraph's content with the MERGEFIELD field-code runs,
# reusing the paragraph's own identity (w14:paraId etc.) and the exact
# run formatting that was on the original run.
$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="49FEE76D" w14:textId="77777777" w:rsidR="00FF1F7F" w:rsidRDefault="00FF1F7F" w:rsidP="00FF1F7F"><w:pPr><w:pStyle w:val="Footer"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>«=website»</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$hitRange.InsertXML($xmlFrag)

Write-Output "Replaced '$targetText' with a =website MERGEFIELD"
